$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial that was bumped by one day
# (45188 -> 45189) for every data row (rows 2 through 379).
$ws.Range("C2:C379").Value = 45189
